$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.949.27"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.375.79"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "3.952.51"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "3.380.87"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "61.043.97"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "3.516.08"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("E27").Value = "  +6.92%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("E33").Value = "  -4.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").Value = "3.410.98"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0762"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "2.435.41"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0259"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("E51").Value = "  +5.49%  "

Write-Host "Updated cryptos list"
